# "Generate Report for Archive" - refresh the localization-status report:
#   1. Status moved from "Ready for handoff" to "In Translation" for the
#      tracked file (Overview sheet's per-language status columns, and the
#      per-language sheets' Status column for the same row).
#   2. The Status column got narrower now that the longer "Ready for
#      handoff" text is gone (columns resized to fit the new text).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- 1. Update the status text everywhere it appears ---------------------
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# --- 2. Resize the Status columns to their new (narrower) width ----------
# Target stored column width is ~13.41 characters; this engine's
# ColumnWidth setter snaps to a 1/6-character pixel grid, so 12.5 is the
# closest achievable input (resolves to ~13.33 stored).
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
